$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -388

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -309

# ALC row 46
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 83333544
$ws.Range("I46").Value = 166666670
$ws.Range("J46").Value = 420
$ws.Range("K46").Value = 500000010
$ws.Range("L46").Value = 1260
$ws.Range("M46").Value = -499999891
$ws.Range("N46").Value = -1498

# ALC row 60
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 83333544
$ws.Range("I60").Value = 166666670
$ws.Range("J60").Value = 420
$ws.Range("K60").Value = 500000010
$ws.Range("L60").Value = 1260
$ws.Range("M60").Value = -499999526
$ws.Range("N60").Value = -2228

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7103
$ws.Range("I132").Value = 6359.727
$ws.Range("J132").Value = 8125
$ws.Range("K132").Value = 19079.181
$ws.Range("L132").Value = 24375
$ws.Range("M132").Value = -16549.181
$ws.Range("N132").Value = -29435

# ARM row 13
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100
$ws.Range("N13").Value = -388

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4663.35
$ws.Range("I32").Value = 4402.1123
$ws.Range("J32").Value = 6777
$ws.Range("K32").Value = 4402.1123
$ws.Range("L32").Value = 6777
$ws.Range("M32").Value = -4115.1123
$ws.Range("N32").Value = -7351

# ARM row 76
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 16644
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 16644
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 16644
$ws.Range("N76").Value = -17320

# ARM row 79
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 16644
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 16644
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 16644
$ws.Range("N79").Value = -18984

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1046.625
$ws.Range("I97").Value = 824.3333
$ws.Range("J97").Value = 1713.5
$ws.Range("K97").Value = 824.3333
$ws.Range("L97").Value = 1713.5
$ws.Range("M97").Value = -328.3333
$ws.Range("N97").Value = -2705.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4285.146
$ws.Range("I132").Value = 1909.4814
$ws.Range("J132").Value = 7339.5713
$ws.Range("K132").Value = 5728.4442
$ws.Range("L132").Value = 22018.7139
$ws.Range("M132").Value = -3198.4442
$ws.Range("N132").Value = -27078.7139

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1702.16
$ws.Range("I107").Value = 1722.5264
$ws.Range("J107").Value = 1637.6666
$ws.Range("K107").Value = 1722.5264
$ws.Range("L107").Value = 1637.6666
$ws.Range("M107").Value = 197.4736
$ws.Range("N107").Value = -5477.6666

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4795.731
$ws.Range("I134").Value = 2642.75
$ws.Range("J134").Value = 6141.3438
$ws.Range("K134").Value = 7928.25
$ws.Range("L134").Value = 18424.0314
$ws.Range("M134").Value = -5393.25
$ws.Range("N134").Value = -23494.0314

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7094142
$ws.Range("I31").Value = 1385.7838
$ws.Range("J31").Value = 33337340
$ws.Range("K31").Value = 1385.7838
$ws.Range("L31").Value = 33337340
$ws.Range("M31").Value = -1090.7838
$ws.Range("N31").Value = -33337930

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7094142
$ws.Range("I34").Value = 1385.7838
$ws.Range("J34").Value = 33337340
$ws.Range("K34").Value = 1385.7838
$ws.Range("L34").Value = 33337340
$ws.Range("M34").Value = -1183.7838
$ws.Range("N34").Value = -33337744

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4074.0527
$ws.Range("I132").Value = 2578.1
$ws.Range("J132").Value = 5736.222
$ws.Range("K132").Value = 7734.299999999999
$ws.Range("L132").Value = 17208.666
$ws.Range("M132").Value = -5204.299999999999
$ws.Range("N132").Value = -22268.666

# CUL row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 83335910
$ws.Range("I25").Value = 1150
$ws.Range("J25").Value = 166670670
$ws.Range("K25").Value = 3450
$ws.Range("L25").Value = 500012010
$ws.Range("M25").Value = -3281
$ws.Range("N25").Value = -500012348

# CUL row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 83335910
$ws.Range("I30").Value = 1150
$ws.Range("J30").Value = 166670670
$ws.Range("K30").Value = 3450
$ws.Range("L30").Value = 500012010
$ws.Range("M30").Value = -3348
$ws.Range("N30").Value = -500012214

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2860.5
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 3350.625
$ws.Range("K132").Value = 8100
$ws.Range("L132").Value = 30155.625
$ws.Range("M132").Value = -5570
$ws.Range("N132").Value = -35215.625

# GSM row 19
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 100000000
$ws.Range("I19").Value = 100000000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 100000000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -99999712

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = ""
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = 0

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2879.0645
$ws.Range("I132").Value = 3451.8
$ws.Range("J132").Value = 2768.923
$ws.Range("K132").Value = 10355.4
$ws.Range("L132").Value = 8306.769
$ws.Range("M132").Value = -7825.400000000001
$ws.Range("N132").Value = -13366.769

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2217.5217
$ws.Range("I82").Value = 1842.8572
$ws.Range("J82").Value = 2800.3333
$ws.Range("K82").Value = 1842.8572
$ws.Range("L82").Value = 2800.3333
$ws.Range("M82").Value = -1481.8572
$ws.Range("N82").Value = -3522.3333

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2217.5217
$ws.Range("I85").Value = 1842.8572
$ws.Range("J85").Value = 2800.3333
$ws.Range("K85").Value = 1842.8572
$ws.Range("L85").Value = 2800.3333
$ws.Range("M85").Value = -594.8571999999999
$ws.Range("N85").Value = -5296.3333

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3649.125
$ws.Range("I100").Value = 3500
$ws.Range("J100").Value = 3698.8333
$ws.Range("K100").Value = 3500
$ws.Range("L100").Value = 3698.8333
$ws.Range("M100").Value = -2959
$ws.Range("N100").Value = -4780.8333

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 41668652
$ws.Range("I136").Value = 55557090
$ws.Range("J136").Value = 3335
$ws.Range("K136").Value = 166671270
$ws.Range("L136").Value = 10005
$ws.Range("M136").Value = -166668720
$ws.Range("N136").Value = -15105

# WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 490
$ws.Range("I14").Value = 512.5
$ws.Range("J14").Value = 475
$ws.Range("K14").Value = 512.5
$ws.Range("L14").Value = 475
$ws.Range("M14").Value = -344.5
$ws.Range("N14").Value = -811

# WVR row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 70011
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 70011
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 70011
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = -70491

# WVR row 30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 31504
$ws.Range("I30").Value = 3750
$ws.Range("J30").Value = 50006.668
$ws.Range("K30").Value = 3750
$ws.Range("L30").Value = 50006.668
$ws.Range("M30").Value = -3643
$ws.Range("N30").Value = -50220.668
